# Añadido funcionalidad de envio por correo
# Appends the newest batch of scraped EVOWHEY PROTEIN price rows (rows 19-37)
# and corrects the timestamp on row 18 (A18) to match the freshly re-scraped value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the previous last row's timestamp (tiny precision correction) ---
$ws.Range("A18").Value = 45818.39371379629

# --- New scraped rows ---
$dates = @(
    45833.41284649305,
    45833.41719120371,
    45833.41743076389,
    45833.41856145833,
    45833.4233575,
    45833.42987278935,
    45833.43047987269,
    45833.43419064815,
    45833.43690237268,
    45833.43811886574,
    45833.43889030092,
    45833.43973851852,
    45833.44137646991,
    45833.44280247685,
    45833.44407189815,
    45833.44527612269,
    45833.45329072917,
    45833.4542096412,
    45833.45888479183
)

$startRow = 19
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = "EVOWHEY PROTEIN"
    $ws.Cells.Item($row, 3).Value = "2Kg"
    $ws.Cells.Item($row, 4).Value = "37,90€"
}
